$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column A header: "data_words" -> "data_words_type"
$ws.Cells.Item(1, 1).Value = "data_words_type"

# Row 2 stays "qa" (already shared-string 11 -> will be reindexed to 10), just
# reassign to normalize / make it explicit.
$ws.Cells.Item(2, 1).Value = "qa"

# Add three new trial rows (3, 4, 5) mirroring the existing param row, varying
# use_bigram (col C) and eta_entry (col G).
$ws.Cells.Item(3, 1).Value = "qa"
$ws.Cells.Item(3, 2).Value = 40
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(3, 4).Value = $false
$ws.Cells.Item(3, 5).Value = 0.1
$ws.Cells.Item(3, 6).Value = 0.25
$ws.Cells.Item(3, 7).Value = "auto"
$ws.Cells.Item(3, 8).Value = 1000
$ws.Cells.Item(3, 9).Value = 2
$ws.Cells.Item(3, 10).Value = 20

$ws.Cells.Item(4, 1).Value = "qa"
$ws.Cells.Item(4, 2).Value = 40
$ws.Cells.Item(4, 3).Value = $false
$ws.Cells.Item(4, 4).Value = $false
$ws.Cells.Item(4, 5).Value = 0.1
$ws.Cells.Item(4, 6).Value = 0.25
$ws.Cells.Item(4, 7).Value = 0.25
$ws.Cells.Item(4, 8).Value = 1000
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(4, 10).Value = 20

$ws.Cells.Item(5, 1).Value = "qa"
$ws.Cells.Item(5, 2).Value = 40
$ws.Cells.Item(5, 3).Value = $true
$ws.Cells.Item(5, 4).Value = $false
$ws.Cells.Item(5, 5).Value = 0.1
$ws.Cells.Item(5, 6).Value = 0.25
$ws.Cells.Item(5, 7).Value = 0.25
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 20

# Match the author's final selection/active cell.
$ws.Range("D7").Select()
